$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their text representation (no auto numeric/date conversion)
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '98.135.04'
$ws.Range("E2").Value = '  +2.97%  '

# Row 3
$ws.Range("D3").Value = '3.619.85'
$ws.Range("E3").Value = '  +1.68%  '

# Row 4
$ws.Range("E4").Value = '  +0.23%  '

# Row 5
$ws.Range("D5").Value = '243.73'
$ws.Range("E5").Value = '  +3.73%  '

# Row 6
$ws.Range("D6").Value = '1.74'
$ws.Range("E6").Value = '  +20.66%  '

# Row 7
$ws.Range("D7").Value = '657.89'
$ws.Range("E7").Value = '  +1.46%  '

# Row 8
$ws.Range("D8").Value = '0.418'
$ws.Range("E8").Value = '  +5.75%  '

# Row 9
$ws.Range("D9").Value = '1.08'
$ws.Range("E9").Value = '  +10.72%  '

# Row 10
$ws.Range("D10").Value = '1.00'
$ws.Range("E10").Value = '  -0.01%  '

# Row 11
$ws.Range("D11").Value = '3.621.60'
$ws.Range("E11").Value = '  +1.75%  '

# Row 12
$ws.Range("D12").Value = '43.99'
$ws.Range("E12").Value = '  +4.93%  '

# Row 13
$ws.Range("D13").Value = '0.206'
$ws.Range("E13").Value = '  +2.64%  '

# Row 14
$ws.Range("D14").Value = '6.47'
$ws.Range("E14").Value = '  -0.47%  '

# Row 15
$ws.Range("D15").Value = '4.297.55'
$ws.Range("E15").Value = '  +1.28%  '

# Row 16
$ws.Range("D16").Value = '98.058.29'
$ws.Range("E16").Value = '  +2.97%  '

# Row 17
$ws.Range("D17").Value = '0.0000260'
$ws.Range("E17").Value = '  +3.74%  '

# Row 18
$ws.Range("D18").Value = '3.622.45'
$ws.Range("E18").Value = '  +1.60%  '

# Row 19
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '12.83'
$ws.Range("E19").Value = '  +1.56%  '

# Row 20
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '7.79'
$ws.Range("E20").Value = '  -0.93%  '

# Row 21
$ws.Range("D21").Value = '18.13'
$ws.Range("E21").Value = '  +3.06%  '

# Row 22
$ws.Range("D22").Value = '0.535'
$ws.Range("E22").Value = '  +14.11%  '

# Row 23
$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D23").Value = '3.47'
$ws.Range("E23").Value = '  +0.46%  '

# Row 24
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '514.62'
$ws.Range("E24").Value = '  +2.34%  '

# Row 25
$ws.Range("D25").Value = '0.0000207'
$ws.Range("E25").Value = '  +8.09%  '

# Row 26
$ws.Range("D26").Value = '6.91'
$ws.Range("E26").Value = '  +5.87%  '

# Row 27
$ws.Range("D27").Value = '99.72'
$ws.Range("E27").Value = '  +9.18%  '

# Row 28
$ws.Range("D28").Value = '12.98'
$ws.Range("E28").Value = '  +4.99%  '

# Row 29
$ws.Range("D29").Value = '3.814.67'
$ws.Range("E29").Value = '  +1.66%  '

# Row 30
$ws.Range("D30").Value = '0.156'
$ws.Range("E30").Value = '  +13.09%  '

# Row 31
$ws.Range("D31").Value = '3.05'
$ws.Range("E31").Value = '  +0.78%  '

# Row 32
$ws.Range("D32").Value = '11.84'
$ws.Range("E32").Value = '  +6.08%  '

# Row 33
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").Value = '  -0.12%  '

# Row 34
$ws.Range("D34").Value = '0.188'
$ws.Range("E34").Value = '  +5.72%  '

# Row 35
$ws.Range("D35").Value = '0.994'
$ws.Range("E35").Value = '  -0.52%  '

# Row 36
$ws.Range("D36").Value = '31.95'
$ws.Range("E36").Value = '  +0.74%  '

# Row 37
$ws.Range("D37").Value = '8.86'
$ws.Range("E37").Value = '  +8.67%  '

# Row 38
$ws.Range("D38").Value = '0.574'
$ws.Range("E38").Value = '  +3.71%  '

# Row 39
$ws.Range("D39").Value = '614.95'
$ws.Range("E39").Value = '  +10.12%  '

# Row 40
$ws.Range("D40").Value = '1.64'
$ws.Range("E40").Value = '  +8.76%  '

# Row 41
$ws.Range("D41").Value = '1.98'
$ws.Range("E41").Value = '  +12.78%  '

# Row 42
$ws.Range("E42").Value = '  +3.14%  '

# Row 43
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  -0.02%  '

# Row 44
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '0.928'
$ws.Range("E44").Value = '  +3.67%  '

# Row 45
$ws.Range("D45").Value = '6.03'
$ws.Range("E45").Value = '  +8.14%  '

# Row 46
$ws.Range("D46").Value = '0.0443'
$ws.Range("E46").Value = '  +9.00%  '

# Row 47
$ws.Range("D47").Value = '2.30'
$ws.Range("E47").Value = '  +1.23%  '

# Row 48
$ws.Range("D48").Value = '23.67'
$ws.Range("E48").Value = '  +0.56%  '

# Row 49
$ws.Range("D49").Value = '8.64'
$ws.Range("E49").Value = '  +7.96%  '

# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.398'
$ws.Range("E50").Value = '  +36.25%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '33.02'
$ws.Range("E51").Value = '  -4.26%  '
